# Generate Report for Handoff
# Rename the generated source/handoff file identifiers and bump their
# timestamps, mirroring what the localization pipeline does when it
# regenerates the handoff artifacts under a new GUID.

$wb = $excel.ActiveWorkbook

$oldGuid = "0e804411-ab76-4683-896f-2a0058e24bde"
$newGuid = "8fe27c90-b8da-483b-bda5-03555ed4961e"

$oldZhCnHash = "0e804411-ab76-4683-896f-2a0058e24bde.80918c25920574b5a82d7f9b9270983c9e6b5a54.zh-cn.xlf"
$newZhCnHash = "8fe27c90-b8da-483b-bda5-03555ed4961e.2cc13e4545c5f9127bd297223bb87d3ef53d1c0e.zh-cn.xlf"

$oldDeDeHash = "0e804411-ab76-4683-896f-2a0058e24bde.80918c25920574b5a82d7f9b9270983c9e6b5a54.de-de.xlf"
$newDeDeHash = "8fe27c90-b8da-483b-bda5-03555ed4961e.2cc13e4545c5f9127bd297223bb87d3ef53d1c0e.de-de.xlf"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("G2").Value = "2016-09-06 10:58:25"

foreach ($hl in $wsOverview.Hyperlinks) {
    $hl.TextToDisplay = "e2e\$newGuid.md"
}

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = $newZhCnHash
$wsZhCn.Range("H2").Value = "2016-09-06 10:58:12"

foreach ($hl in $wsZhCn.Hyperlinks) {
    $hl.TextToDisplay = "$newGuid.md"
}

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = $newDeDeHash
# H2 shares its text with Overview!G2 ("Latest HO Xliff Generate Date" /
# "Latest Handoff Datetime"), so it is already updated above.

foreach ($hl in $wsDeDe.Hyperlinks) {
    $hl.TextToDisplay = "$newGuid.md"
}
